$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
$ws.Range("C2").Value = 45174
$ws.Range("C3").Value = 45174
$ws.Range("C4").Value = 45174
$ws.Range("C5").Value = 45174
